$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.3639155626297
$ws.Range("B1").Value = 3.027185440063477
$ws.Range("C1").Value = 3.097575187683105
$ws.Range("D1").Value = 1.096873044967651
$ws.Range("E1").Value = 0.8176232576370239
